# Actualización automática 2025-07-24 16:45:08
$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("M5").Value = 4521.91
$ws1.Range("C25").Value = 1534.46
$ws1.Range("M32").Value = 3369.61
$ws1.Range("C36").Value = 1036.8
$ws1.Range("C56").Value = "4 de 54"

# --- Sheet: VENTA MENSUAL ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F5").Value = 7353.3
$ws2.Range("F25").Value = 11643.33
$ws2.Range("F32").Value = 3369.61
$ws2.Range("F36").Value = 14420.72
$ws2.Range("F56").Value = 62535.26

# --- Sheet: CUMPLIMIENTO MENSUAL ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D2").Value = 4644.86
$ws3.Range("E2").Value = 5325.48304517915
$ws3.Range("F2").Value = 0.4658676215003332

$ws3.Range("D16").Value = 50722.88
$ws3.Range("E16").Value = 1103.580000000002
$ws3.Range("F16").Value = 0.9787062438761975

$ws3.Range("D19").Value = 62535.25999999999
$ws3.Range("E19").Value = 51171.19064517916
$ws3.Range("F19").Value = 0.5499710847112905
